$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 15 (pushes the existing SUM row down to row 16) ---
$ws.Rows("15").Insert()

# --- New row 15: Foam board ---
$ws.Range("A15").Value = "Foam board"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 1.5
$ws.Range("C15").NumberFormat = """$""#,##0.00"
$ws.Range("D15").Formula = "=B15*C15"
$ws.Range("D15").NumberFormat = """$""#,##0.00"

# --- Fix the SUM formula on (now) row 16 to include the new row ---
$ws.Range("D16").Formula = "=SUM(D2:D15)"

# --- Update item descriptions that became more detailed ---
$ws.Range("A12").Value = "9-DOF board (Gyro, accel, mag)"
$ws.Range("A13").Value = "Feather M0 LoRa (Arduino)"
$ws.Range("A14").Value = "RFM95W (Radio)"

# --- New column E: "COST IS IN USD" header ---
$ws.Range("E1").Value = "COST IS IN USD"

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 26.75
$ws.Columns("E").ColumnWidth = 12.75

# --- Selection moves to A19 ---
$ws.Range("A19").Select()
